# Check the four to-do checkboxes (content-control checkboxes) that were
# marked complete in this revision:
#   1524207776 - "Remove admin from the SDS."
#   1822843809 - "Update architectural diagram and the diagrams in general."
#   1361858848 - "Also, to open the project you do 'python application' ..."
#   642702992  - "Add range for avg glucose level between 55.1 and 272"
#
# Checking a checkbox content control requires two things in the OOXML:
#   1. w14:checked/@w14:val flips from 0 to 1 (the logical state)
#   2. the displayed glyph run flips from ☐ (U+2610) to ☒ (U+2612)
# Setting ContentControl.Checked handles (1); we also update the glyph
# text ourselves so the visible symbol matches the new state.

$d = $word.ActiveDocument

$idsToCheck = @("1524207776", "1822843809", "1361858848", "642702992")

$checkedGlyph = [char]0x2612

$ccs = $d.ContentControls
for ($i = 1; $i -le $ccs.Count; $i++) {
    $cc = $ccs.Item($i)
    if ($idsToCheck -contains "$($cc.ID)") {
        $cc.Checked = $true
        $cc.Range.Text = $checkedGlyph
    }
}
